$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 0.0002777777777777778
$ws.Range("K2").Value = 1888
$ws.Range("L2").Value = 0.003776
